# Daily attendance processing - 2025-12-01 20:28:06
# Reverses the order of the comma-separated "Recorded By" entries in
# column G (e.g. "System, user@example.com" -> "user@example.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $joined = ""
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                if ($joined -ne "") { $joined += ", " }
                $joined += $parts[$i]
            }
            $cell.Value2 = $joined
        }
    }
}
